# Update Name of Algo
# Apply updated values to result_data_KNN worksheet (imputed values changed
# after updating the algorithm name / re-running imputation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.228
$ws.Range("E4").Value = 13.265
$ws.Range("A9").Value = -20.912
$ws.Range("E10").Value = 12.633
$ws.Range("A18").Value = -21.825
$ws.Range("A20").Value = -21.828
$ws.Range("D21").Value = -7.675999999999999
